# --------------------------------------------------------------------------
# Updates the "two-digit divided by one-digit" worksheet:
#   1. Bumps the heading date from 2025-07-15 Tuesday -> 2025-07-16 Wednesday.
#   2. Replaces every division-problem answer in the 5x5 data table with the
#      newly generated value for that same (row, column) slot.
#
# Cells are addressed by (row, column) via Table.Cell(r, c).Range.Text rather
# than by a document-wide text Find & Replace, because several of the new
# answer strings coincide with *other* cells' old answer strings (e.g. the
# text "46÷6=7, 4" is both an old value at one cell and a new value at a
# different cell). A sequential find/replace-all over the whole document
# would clobber the wrong cell in that situation, so direct cell addressing
# is used to guarantee each slot gets exactly the value it should.
# --------------------------------------------------------------------------

$d = $word.ActiveDocument

# 1. Heading date.
$d.Content.Find.Execute("2025-07-15 Tuesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-07-16 Wednesday", 2) | Out-Null

# 2. Table answers. Only every 4th row (1, 5, 9, 13, 17) holds data; the rows
#    in between are blank spacer rows and are left untouched.
$t = $d.Tables.Item(1)

function Set-Answer($row, $col, $oldText, $newText) {
    $cell = $t.Cell($row, $col)
    $current = $cell.Range.Text
    if ($current -notlike "$oldText*") {
        Write-Host "WARNING: Cell($row,$col) expected '$oldText' but found '$current'"
    }
    $cell.Range.Text = $newText
}

Set-Answer 1  1 "16÷9=1, 7"   "46÷6=7, 4"
Set-Answer 1  2 "46÷4=11, 2"  "45÷3=15, 0"
Set-Answer 1  3 "51÷2=25, 1"  "39÷7=5, 4"
Set-Answer 1  4 "48÷6=8, 0"   "15÷5=3, 0"
Set-Answer 1  5 "52÷3=17, 1"  "46÷4=11, 2"

Set-Answer 5  1 "88÷7=12, 4"  "22÷8=2, 6"
Set-Answer 5  2 "94÷2=47, 0"  "87÷9=9, 6"
Set-Answer 5  3 "22÷6=3, 4"   "80÷2=40, 0"
Set-Answer 5  4 "66÷5=13, 1"  "81÷8=10, 1"
Set-Answer 5  5 "46÷6=7, 4"   "46÷8=5, 6"

Set-Answer 9  1 "27÷4=6, 3"   "36÷6=6, 0"
Set-Answer 9  2 "97÷8=12, 1"  "74÷7=10, 4"
Set-Answer 9  3 "43÷2=21, 1"  "31÷9=3, 4"
Set-Answer 9  4 "90÷3=30, 0"  "78÷3=26, 0"
Set-Answer 9  5 "77÷6=12, 5"  "31÷5=6, 1"

Set-Answer 13 1 "24÷5=4, 4"   "70÷4=17, 2"
Set-Answer 13 2 "50÷7=7, 1"   "90÷6=15, 0"
Set-Answer 13 3 "71÷6=11, 5"  "68÷9=7, 5"
Set-Answer 13 4 "22÷3=7, 1"   "37÷8=4, 5"
Set-Answer 13 5 "39÷4=9, 3"   "37÷5=7, 2"

Set-Answer 17 1 "81÷4=20, 1"  "67÷2=33, 1"
Set-Answer 17 2 "79÷9=8, 7"   "70÷2=35, 0"
Set-Answer 17 3 "72÷6=12, 0"  "68÷2=34, 0"
Set-Answer 17 4 "93÷7=13, 2"  "44÷6=7, 2"
Set-Answer 17 5 "38÷7=5, 3"   "59÷6=9, 5"

Write-Host "Edit complete"
